$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 147, shifting rows 147:188 down to 148:189.
$ws.Rows.Item(147).Insert()

# Populate the new row 147 with a copy of the (now shifted-down) row 148,
# except for the new date and the updated price figures.
# (Value2 is used for the reads -- Value's getter is unreliable in this host.)
$ws.Cells.Item(147, 1).Value = $ws.Cells.Item(148, 1).Value2   # Mercado ID
$ws.Cells.Item(147, 2).Value = $ws.Cells.Item(148, 2).Value2   # Mercado
$ws.Cells.Item(147, 3).Value = $ws.Cells.Item(148, 3).Value2   # Region
$ws.Cells.Item(147, 4).Value = 44559                           # Fecha
$ws.Cells.Item(147, 5).Value = $ws.Cells.Item(148, 5).Value2   # Codreg
$ws.Cells.Item(147, 6).Value = $ws.Cells.Item(148, 6).Value2   # Categoria ID
$ws.Cells.Item(147, 7).Value = $ws.Cells.Item(148, 7).Value2   # Categoria
$ws.Cells.Item(147, 8).Value = $ws.Cells.Item(148, 8).Value2   # Variedad
$ws.Cells.Item(147, 9).Value = $ws.Cells.Item(148, 9).Value2   # Calidad
$ws.Cells.Item(147, 10).Value = $ws.Cells.Item(148, 10).Value2 # Volumen
$ws.Cells.Item(147, 11).Value = 800                             # Precio minimo
$ws.Cells.Item(147, 12).Value = 800                             # Precio maximo
$ws.Cells.Item(147, 13).Value = 800                             # Precio promedio ponderado
$ws.Cells.Item(147, 14).Value = $ws.Cells.Item(148, 14).Value2 # Unidad de comercializacion
$ws.Cells.Item(147, 15).Value = $ws.Cells.Item(148, 15).Value2 # Origen
$ws.Cells.Item(147, 16).Value = 800                             # Precio $/Kg
$ws.Cells.Item(147, 17).Value = $ws.Cells.Item(148, 17).Value2 # Kg o Unidades
$ws.Cells.Item(147, 18).Value = $ws.Cells.Item(148, 18).Value2 # Clasificacion
